$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new daily electricity price data (automatic update).
$ws.Range("A2").Value = 45900
$ws.Range("B2").Value = 81.06
$ws.Range("C2").Value = 75.40000000000001
$ws.Range("D2").Value = 70.5
$ws.Range("E2").Value = 67
$ws.Range("F2").Value = 70.12
$ws.Range("G2").Value = 70.12
$ws.Range("H2").Value = 70.12
$ws.Range("I2").Value = 75.40000000000001
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 3.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -0.01
$ws.Range("N2").Value = -0.08
$ws.Range("O2").Value = -0.95
$ws.Range("P2").Value = -0.6
$ws.Range("Q2").Value = -0.62
$ws.Range("R2").Value = -0.98
$ws.Range("S2").Value = -0.27
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 3.52
$ws.Range("V2").Value = 59.55
$ws.Range("W2").Value = 73.88
$ws.Range("X2").Value = 71.40000000000001
$ws.Range("Y2").Value = 62
$ws.Range("Z2").Value = 37.3
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 73.48999999999999
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 78.23
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 72.76000000000001
$ws.Range("AG2").Value = "9h-19h"
